# ngClick & ngDblClick Directives
# Splits the long "Directives" cell run so that "7.ng-click      8.ng-dblclick"
# becomes its own green (00B050) run, matching the other directive "headings"
# in that cell, while the surrounding text keeps its original black color.

$d = $word.ActiveDocument

$fullText = "          7.ng-click      8.ng-dblclick               9.ng-if                          10.ng-show  11.ng-hide                  12.ng-switch          13.ng-submit  14.ng-include        15.ng-cloak            16.ng-mouseover                                                                  17.ng-mouseleave                                                            18.ng-mouseenter   "
$part2 = "7.ng-click      8.ng-dblclick"

# Locate the run's full text. MatchWholeWord must be $false: the text starts
# with whitespace, and this Find engine fails to match whole-word text that
# begins/ends on whitespace when MatchWholeWord is requested.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute($fullText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target Directives text"
}

$s = $rng.Start
$e = $rng.End

# Force a genuine content change (and restore) over the whole matched range
# first. This drops the legacy w:rsidRPr carried by the original run so the
# freshly split-out runs come out as plain <w:r> (no rsid), matching a real
# Word edit.
$placeholder = ""
for ($i = 0; $i -lt $fullText.Length; $i++) { $placeholder += "*" }
$rng.Text = $placeholder
$rng2 = $d.Range($s, $s + $fullText.Length)
$rng2.Text = $fullText

# Re-locate the "7.ng-click      8.ng-dblclick" substring inside the restored
# range and recolor it green (RGB 00B050 -> Word BGR long 5287936), which
# splits it off into its own run while leaving the rest of the text (still
# black / themeColor text1) as the surrounding run(s).
$p2Start = $s + $fullText.IndexOf($part2)
$p2End = $p2Start + $part2.Length
$colorRng = $d.Range($p2Start, $p2End)
$colorRng.Font.Color = 5287936
